# Update "想去人数" (column F) counts across all sheets to reflect the
# newly generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 1209
        4  = 50
        5  = 1352
        6  = 1707
        7  = 6227
        9  = 1824
        10 = 485
        11 = 6
        12 = 18
        15 = 26
        16 = 6954
        17 = 127
        18 = 55
        19 = 168
        20 = 104
        21 = 1711
        22 = 841
        23 = 16
        24 = 44
        25 = 163
        26 = 1589
        27 = 763
        28 = 317
        33 = 3896
    }
    "演出" = @{
        8 = 444
    }
    "本地生活" = @{
        3 = 2261
        4 = 661
    }
    "全部类型" = @{
        3  = 2261
        4  = 661
        5  = 1209
        7  = 50
        10 = 1352
        12 = 1707
        13 = 6227
        14 = 1824
        17 = 485
        19 = 18
        23 = 6954
        24 = 127
        25 = 55
        26 = 168
        27 = 104
        28 = 1711
        29 = 841
        30 = 16
        31 = 44
        32 = 163
        33 = 1589
        34 = 763
        36 = 317
        44 = 3896
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
